$d = $word.ActiveDocument

# The diff only reorders the <w:rPr> child elements inside a handful of
# character styles (styles.xml) so they follow the wml.xsd sequence
# (w:b / w:i before w:color, etc.) -- OOXMLValidator flagged the old
# "color before b/i" order. No formatting value actually changes;
# re-assigning each style's own (already-true) Font.Bold / Font.Italic
# value is a no-op semantically, but forces Word to re-emit that
# style's rPr in canonical, schema-valid element order on save.

$boldOnly = @(
    "KeywordTok",
    "ImportTok",
    "ControlFlowTok",
    "AlertTok",
    "ErrorTok"
)

$italicOnly = @(
    "CommentTok",
    "DocumentationTok"
)

$boldAndItalic = @(
    "AnnotationTok",
    "CommentVarTok",
    "InformationTok",
    "WarningTok"
)

foreach ($id in $boldOnly) {
    $s = $d.Styles($id)
    $s.Font.Bold = $s.Font.Bold
}

foreach ($id in $italicOnly) {
    $s = $d.Styles($id)
    $s.Font.Italic = $s.Font.Italic
}

foreach ($id in $boldAndItalic) {
    $s = $d.Styles($id)
    $s.Font.Bold = $s.Font.Bold
    $s.Font.Italic = $s.Font.Italic
}
